$d = $word.ActiveDocument

$old = "{% if inputDateHeader %}{{ inputDateHeader }}{% if paperInputDateHeader %}({{ paperInputDateHeader }}){% else %}{% endif %}{% else %}-{% endif %}"
$new = '{% if inputDateHeader %}{{ inputDateHeader | date("dd.MM.YYYY") }}{% if paperInputDateHeader %} ({{ paperInputDateHeader | date("dd.MM.YYYY") }}){% else %}{% endif %}{% else %}-{% endif %}'

$range = $d.Content
$found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $range.Text = $new
    Write-Output "Replaced inputDateHeader block successfully."
} else {
    Write-Output "WARNING: target text not found!"
}
